# Update loading_percent values for Case_2_131 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 14.38867885856291
$ws.Cells.Item(2, 3).Value = 13.35300972362557
$ws.Cells.Item(2, 5).Value = 17.32057726772298
$ws.Cells.Item(2, 6).Value = 35.81139896153815
$ws.Cells.Item(2, 7).Value = 27.443320671177
$ws.Cells.Item(2, 8).Value = 13.96263339366994
$ws.Cells.Item(2, 10).Value = 7.510709608096124
$ws.Cells.Item(2, 12).Value = 12.66341527937687
$ws.Cells.Item(2, 13).Value = 16.19574036529395
$ws.Cells.Item(2, 14).Value = 18.17675094538633
$ws.Cells.Item(2, 15).Value = 21.07757178974234
$ws.Cells.Item(3, 2).Value = 13.96629500281782
$ws.Cells.Item(3, 3).Value = 13.31289895905496
$ws.Cells.Item(3, 5).Value = 17.36275781979108
$ws.Cells.Item(3, 6).Value = 35.85350203704393
$ws.Cells.Item(3, 7).Value = 27.46509393264731
$ws.Cells.Item(3, 8).Value = 14.00474043387731
$ws.Cells.Item(3, 10).Value = 7.497327120999556
$ws.Cells.Item(3, 12).Value = 12.65978683440345
$ws.Cells.Item(3, 13).Value = 16.1047820188265
$ws.Cells.Item(3, 14).Value = 18.22244835448211
$ws.Cells.Item(3, 15).Value = 21.1377559453253
$ws.Cells.Item(4, 2).Value = 13.70194218217979
$ws.Cells.Item(4, 3).Value = 13.2885408467938
$ws.Cells.Item(4, 5).Value = 17.39069504065285
$ws.Cells.Item(4, 6).Value = 35.88756571033765
$ws.Cells.Item(4, 7).Value = 27.48761922331633
$ws.Cells.Item(4, 8).Value = 14.03285630918482
$ws.Cells.Item(4, 10).Value = 7.489025356693467
$ws.Cells.Item(4, 12).Value = 12.65901708404278
$ws.Cells.Item(4, 13).Value = 16.0506510436668
$ws.Cells.Item(4, 14).Value = 18.25223925446672
$ws.Cells.Item(4, 15).Value = 21.17935611996302
$ws.Cells.Item(5, 2).Value = 13.59311733966401
$ws.Cells.Item(5, 3).Value = 13.27868650479176
$ws.Cells.Item(5, 5).Value = 17.40259278979891
$ws.Cells.Item(5, 6).Value = 35.90351022332381
$ws.Cells.Item(5, 7).Value = 27.49909552103829
$ws.Cells.Item(5, 8).Value = 14.04488238661984
$ws.Cells.Item(5, 10).Value = 7.485621451177679
$ws.Cells.Item(5, 12).Value = 12.65907159794855
$ws.Cells.Item(5, 13).Value = 16.02904035816973
$ws.Cells.Item(5, 14).Value = 18.26481578044961
$ws.Cells.Item(5, 15).Value = 21.19747464907024
$ws.Cells.Item(6, 2).Value = 13.57498558448069
$ws.Cells.Item(6, 3).Value = 13.27705464657734
$ws.Cells.Item(6, 5).Value = 17.40459940940721
$ws.Cells.Item(6, 6).Value = 35.90628234916137
$ws.Cells.Item(6, 7).Value = 27.50113969894203
$ws.Cells.Item(6, 8).Value = 14.04691364651381
$ws.Cells.Item(6, 10).Value = 7.485054988435034
$ws.Cells.Item(6, 12).Value = 12.65910293060307
$ws.Cells.Item(6, 13).Value = 16.02547944684833
$ws.Cells.Item(6, 14).Value = 18.26693049099638
$ws.Cells.Item(6, 15).Value = 21.20055357555364
$ws.Cells.Item(7, 2).Value = 13.70047876194002
$ws.Cells.Item(7, 3).Value = 13.28840765195526
$ws.Cells.Item(7, 5).Value = 17.39085341953546
$ws.Cells.Item(7, 6).Value = 35.88777239230998
$ws.Cells.Item(7, 7).Value = 27.48776470431732
$ws.Cells.Item(7, 8).Value = 14.03301619493715
$ws.Cells.Item(7, 10).Value = 7.488979534317481
$ws.Cells.Item(7, 12).Value = 12.65901632648724
$ws.Cells.Item(7, 13).Value = 16.05035775904469
$ws.Cells.Item(7, 14).Value = 18.25240709743167
$ws.Cells.Item(7, 15).Value = 21.17959575445948
$ws.Cells.Item(8, 2).Value = 14.2441728692459
$ws.Cells.Item(8, 3).Value = 13.33912582105101
$ws.Cells.Item(8, 5).Value = 17.33469844866497
$ws.Cells.Item(8, 6).Value = 35.82421083547555
$ws.Cells.Item(8, 7).Value = 27.44892491508554
$ws.Cells.Item(8, 8).Value = 13.97668233210154
$ws.Cells.Item(8, 10).Value = 7.506113221421179
$ws.Cells.Item(8, 12).Value = 12.66186250747
$ws.Cells.Item(8, 13).Value = 16.16403126446057
$ws.Cells.Item(8, 14).Value = 18.1921483336444
$ws.Cells.Item(8, 15).Value = 21.09735745024443
$ws.Cells.Item(9, 2).Value = 15.2641593107929
$ws.Cells.Item(9, 3).Value = 13.44055625217452
$ws.Cells.Item(9, 5).Value = 17.24072467835708
$ws.Cells.Item(9, 6).Value = 35.76477478599995
$ws.Cells.Item(9, 7).Value = 27.44558006574687
$ws.Cells.Item(9, 8).Value = 13.8841666842063
$ws.Cells.Item(9, 10).Value = 7.53901910951183
$ws.Cells.Item(9, 12).Value = 12.678942860434
$ws.Cells.Item(9, 13).Value = 16.39984053492392
$ws.Cells.Item(9, 14).Value = 18.08769003163704
$ws.Cells.Item(9, 15).Value = 20.973059576474
$ws.Cells.Item(10, 2).Value = 15.97754604730559
$ws.Cells.Item(10, 3).Value = 13.51605074482988
$ws.Cells.Item(10, 5).Value = 17.1814892415271
$ws.Cells.Item(10, 6).Value = 35.76088027171566
$ws.Cells.Item(10, 7).Value = 27.48763595550621
$ws.Cells.Item(10, 8).Value = 13.82715046135472
$ws.Cells.Item(10, 10).Value = 7.562746967213966
$ws.Cells.Item(10, 12).Value = 12.69840246571009
$ws.Cells.Item(10, 13).Value = 16.57993681525885
$ws.Cells.Item(10, 14).Value = 18.0192499487496
$ws.Cells.Item(10, 15).Value = 20.90440542734165
$ws.Cells.Item(11, 2).Value = 16.29284410711431
$ws.Cells.Item(11, 3).Value = 13.55055683733536
$ws.Cells.Item(11, 5).Value = 17.15666338692122
$ws.Cells.Item(11, 6).Value = 35.76773185352226
$ws.Cells.Item(11, 7).Value = 27.51642247962797
$ws.Cells.Item(11, 8).Value = 13.80359289861052
$ws.Cells.Item(11, 10).Value = 7.57343859146651
$ws.Cells.Item(11, 12).Value = 12.7087322607727
$ws.Cells.Item(11, 13).Value = 16.66312964091167
$ws.Cells.Item(11, 14).Value = 17.98990746909915
$ws.Cells.Item(11, 15).Value = 20.87811893861292
$ws.Cells.Item(12, 2).Value = 16.41080923968893
$ws.Cells.Item(12, 3).Value = 13.56364251584231
$ws.Cells.Item(12, 5).Value = 17.14756683959048
$ws.Cells.Item(12, 6).Value = 35.77156366749431
$ws.Cells.Item(12, 7).Value = 27.52870753808123
$ws.Cells.Item(12, 8).Value = 13.79501457570224
$ws.Cells.Item(12, 10).Value = 7.577472006410003
$ws.Cells.Item(12, 12).Value = 12.71285413558557
$ws.Cells.Item(12, 13).Value = 16.69479523381878
$ws.Cells.Item(12, 14).Value = 17.97905303392034
$ws.Cells.Item(12, 15).Value = 20.86887760864871
$ws.Cells.Item(13, 2).Value = 16.38546855521387
$ws.Cells.Item(13, 3).Value = 13.56082350935808
$ws.Cells.Item(13, 5).Value = 17.14951241033734
$ws.Cells.Item(13, 6).Value = 35.77068344343705
$ws.Cells.Item(13, 7).Value = 27.52600025336062
$ws.Cells.Item(13, 8).Value = 13.79684683807781
$ws.Cells.Item(13, 10).Value = 7.576604026525464
$ws.Cells.Item(13, 12).Value = 12.71195710440396
$ws.Cells.Item(13, 13).Value = 16.6879685819417
$ws.Cells.Item(13, 14).Value = 17.98137931418255
$ws.Cells.Item(13, 15).Value = 20.87083616698895
$ws.Cells.Item(14, 2).Value = 16.30257841937852
$ws.Cells.Item(14, 3).Value = 13.55163304975909
$ws.Cells.Item(14, 5).Value = 17.15590890932749
$ws.Cells.Item(14, 6).Value = 35.76802232144605
$ws.Cells.Item(14, 7).Value = 27.51740546734215
$ws.Cells.Item(14, 8).Value = 13.80288029089093
$ws.Cells.Item(14, 10).Value = 7.573770732234861
$ws.Cells.Item(14, 12).Value = 12.70906717043152
$ws.Cells.Item(14, 13).Value = 16.66573165017689
$ws.Cells.Item(14, 14).Value = 17.98900932271426
$ws.Cells.Item(14, 15).Value = 20.87734435407784
$ws.Cells.Item(15, 2).Value = 16.25161647903856
$ws.Cells.Item(15, 3).Value = 13.54600596385991
$ws.Cells.Item(15, 5).Value = 17.15986658567713
$ws.Cells.Item(15, 6).Value = 35.76655333915604
$ws.Cells.Item(15, 7).Value = 27.51232102427014
$ws.Cells.Item(15, 8).Value = 13.80662055489264
$ws.Cells.Item(15, 10).Value = 7.572033252157094
$ws.Cells.Item(15, 12).Value = 12.70732430828381
$ws.Cells.Item(15, 13).Value = 16.65213142936586
$ws.Cells.Item(15, 14).Value = 17.99371636333081
$ws.Cells.Item(15, 15).Value = 20.88142368188128
$ws.Cells.Item(16, 2).Value = 15.95674573593061
$ws.Cells.Item(16, 3).Value = 13.51379859157555
$ws.Cells.Item(16, 5).Value = 17.18315430121409
$ws.Cells.Item(16, 6).Value = 35.76060581288294
$ws.Cells.Item(16, 7).Value = 27.48594866051286
$ws.Cells.Item(16, 8).Value = 13.82873790571613
$ws.Cells.Item(16, 10).Value = 7.562046138152676
$ws.Cells.Item(16, 12).Value = 12.69775693035331
$ws.Cells.Item(16, 13).Value = 16.57452375941626
$ws.Cells.Item(16, 14).Value = 18.02120353702351
$ws.Cells.Item(16, 15).Value = 20.90622296517596
$ws.Cells.Item(17, 2).Value = 15.77341126834829
$ws.Cells.Item(17, 3).Value = 13.4940790356729
$ws.Cells.Item(17, 5).Value = 17.1979833664694
$ws.Cells.Item(17, 6).Value = 35.75916415557869
$ws.Cells.Item(17, 7).Value = 27.4722406034881
$ws.Cells.Item(17, 8).Value = 13.84291579609636
$ws.Cells.Item(17, 10).Value = 7.555892769543503
$ws.Cells.Item(17, 12).Value = 12.69226451269285
$ws.Cells.Item(17, 13).Value = 16.52722424514687
$ws.Cells.Item(17, 14).Value = 18.03852433953284
$ws.Cells.Item(17, 15).Value = 20.92270428048616
$ws.Cells.Item(18, 2).Value = 15.6670993794447
$ws.Cells.Item(18, 3).Value = 13.48275250892499
$ws.Cells.Item(18, 5).Value = 17.20671228641539
$ws.Cells.Item(18, 6).Value = 35.75914683278504
$ws.Cells.Item(18, 7).Value = 27.4652652427387
$ws.Cells.Item(18, 8).Value = 13.85129448977729
$ws.Cells.Item(18, 10).Value = 7.55234391995963
$ws.Cells.Item(18, 12).Value = 12.68924463363759
$ws.Cells.Item(18, 13).Value = 16.50013919013833
$ws.Cells.Item(18, 14).Value = 18.04865544666364
$ws.Cells.Item(18, 15).Value = 20.93264921677556
$ws.Cells.Item(19, 2).Value = 15.63095919379366
$ws.Cells.Item(19, 3).Value = 13.47892036444657
$ws.Cells.Item(19, 5).Value = 17.20970205063876
$ws.Cells.Item(19, 6).Value = 35.75928047712195
$ws.Cells.Item(19, 7).Value = 27.46305974912855
$ws.Cells.Item(19, 8).Value = 13.8541698262885
$ws.Cells.Item(19, 10).Value = 7.551140709545955
$ws.Cells.Item(19, 12).Value = 12.68824613158394
$ws.Cells.Item(19, 13).Value = 16.49098992962109
$ws.Cells.Item(19, 14).Value = 18.05211464976271
$ws.Cells.Item(19, 15).Value = 20.93609626329745
$ws.Cells.Item(20, 2).Value = 15.79301762057138
$ws.Cells.Item(20, 3).Value = 13.49617662987566
$ws.Cells.Item(20, 5).Value = 17.19638412928633
$ws.Cells.Item(20, 6).Value = 35.75923361278766
$ws.Cells.Item(20, 7).Value = 27.4736057800181
$ws.Cells.Item(20, 8).Value = 13.84138335594033
$ws.Cells.Item(20, 10).Value = 7.55654880319275
$ws.Cells.Item(20, 12).Value = 12.69283479776727
$ws.Cells.Item(20, 13).Value = 16.53224704352389
$ws.Cells.Item(20, 14).Value = 18.03666306176902
$ws.Cells.Item(20, 15).Value = 20.92090164405651
$ws.Cells.Item(21, 2).Value = 16.32696489173724
$ws.Cells.Item(21, 3).Value = 13.55433203180982
$ws.Cells.Item(21, 5).Value = 17.15402184327126
$ws.Cells.Item(21, 6).Value = 35.76877040479773
$ws.Cells.Item(21, 7).Value = 27.51989243931647
$ws.Cells.Item(21, 8).Value = 13.80109882451951
$ws.Cells.Item(21, 10).Value = 7.574603358405457
$ws.Cells.Item(21, 12).Value = 12.70991032777705
$ws.Cells.Item(21, 13).Value = 16.67225893074181
$ws.Cells.Item(21, 14).Value = 17.98676123651694
$ws.Cells.Item(21, 15).Value = 20.87541338380105
$ws.Cells.Item(22, 2).Value = 16.66754581092109
$ws.Cells.Item(22, 3).Value = 13.59244840217404
$ws.Cells.Item(22, 5).Value = 17.12811007075271
$ws.Cells.Item(22, 6).Value = 35.78221307100858
$ws.Cells.Item(22, 7).Value = 27.55820852251958
$ws.Cells.Item(22, 8).Value = 13.77676651869711
$ws.Cells.Item(22, 10).Value = 7.586313777172705
$ws.Cells.Item(22, 12).Value = 12.72229428058144
$ws.Cells.Item(22, 13).Value = 16.76470134265583
$ws.Cells.Item(22, 14).Value = 17.95564477343577
$ws.Cells.Item(22, 15).Value = 20.84983954245055
$ws.Cells.Item(23, 2).Value = 16.48657044025946
$ws.Cells.Item(23, 3).Value = 13.57209655416985
$ws.Cells.Item(23, 5).Value = 17.14177747039585
$ws.Cells.Item(23, 6).Value = 35.77437991298174
$ws.Cells.Item(23, 7).Value = 27.53702241327816
$ws.Cells.Item(23, 8).Value = 13.7895704159622
$ws.Cells.Item(23, 10).Value = 7.580072062178875
$ws.Cells.Item(23, 12).Value = 12.71557350350353
$ws.Cells.Item(23, 13).Value = 16.71528398318455
$ws.Cells.Item(23, 14).Value = 17.97211542750418
$ws.Cells.Item(23, 15).Value = 20.86310804040813
$ws.Cells.Item(24, 2).Value = 15.78415641652727
$ws.Cells.Item(24, 3).Value = 13.49522827446583
$ws.Cells.Item(24, 5).Value = 17.19710651052629
$ws.Cells.Item(24, 6).Value = 35.75919968322233
$ws.Cells.Item(24, 7).Value = 27.47298576259109
$ws.Cells.Item(24, 8).Value = 13.84207546307868
$ws.Cells.Item(24, 10).Value = 7.556252245209308
$ws.Cells.Item(24, 12).Value = 12.69257654253303
$ws.Cells.Item(24, 13).Value = 16.52997589825435
$ws.Cells.Item(24, 14).Value = 18.03750400605029
$ws.Cells.Item(24, 15).Value = 20.92171515331773
$ws.Cells.Item(25, 2).Value = 14.99402607784165
$ws.Cells.Item(25, 3).Value = 13.41292886530327
$ws.Cells.Item(25, 5).Value = 17.26442224593746
$ws.Cells.Item(25, 6).Value = 35.77386731021317
$ws.Cells.Item(25, 7).Value = 27.43866960807112
$ws.Cells.Item(25, 8).Value = 13.90727142272108
$ws.Cells.Item(25, 10).Value = 7.530194010368477
$ws.Cells.Item(25, 12).Value = 12.67310091788895
$ws.Cells.Item(25, 13).Value = 16.33477023785949
$ws.Cells.Item(25, 14).Value = 18.08769003163704
$ws.Cells.Item(25, 15).Value = 20.973059576474
